# Applies the "[base commands] - [assertMatch(text,regex)]" change to the
# jmeter-showcase.xlsx "#system" reference sheet:
#   1. Removes the obsolete "tn.5250" category column (was column AA),
#      which shifts web/webalert/webcookie/ws/ws.async/xml one column left.
#   2. Removes the now orphaned "tn.5250" entry from the "target" category
#      list in column A.
#   3. Adds the new "assertMatch(text,regex)" command to the "base" category
#      (column F), keeping the alphabetical ordering.
#   4. Adds the new "openFile(filePath)" command to the "external" category
#      (column J), keeping the alphabetical ordering.
#   5. Re-points the affected named ranges to match the new data extents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

function Set-ColumnValues($Sheet, $Column, $StartRow, $Values) {
    $count = $Values.Length
    $endRow = $StartRow + $count - 1
    $arr = New-Object 'object[,]' $count,1
    for ($i = 0; $i -lt $count; $i++) {
        $arr[$i,0] = $Values[$i]
    }
    $addr = "$Column$StartRow" + ":" + "$Column$endRow"
    $range = $Sheet.Range($addr)
    $range.Value2 = $arr
}

# ---------------------------------------------------------------------------
# 1. Delete the "tn.5250" column (column AA). This shifts web, webalert,
#    webcookie, ws, ws.async and xml one column to the left (AB->AA, ... ,
#    AG->AF) automatically.
# ---------------------------------------------------------------------------
$ws.Columns("AA").Delete()

# ---------------------------------------------------------------------------
# 2. Rewrite the "target" category list (column A) without "tn.5250".
#    It previously spanned A2:A33 (32 entries); it now spans A2:A32 (31).
# ---------------------------------------------------------------------------
$targetList = @(
    "aws.s3", "aws.ses", "aws.sqs", "aws.vision", "base", "csv", "desktop",
    "excel", "external", "image", "io", "jms", "json", "localdb", "macro",
    "mail", "number", "pdf", "rdbms", "redis", "sms", "sound", "ssh", "step",
    "step.inTime", "web", "webalert", "webcookie", "ws", "ws.async", "xml"
)
Set-ColumnValues $ws "A" 2 $targetList
$ws.Range("A33").ClearContents()

# ---------------------------------------------------------------------------
# 3. Rewrite the "base" category list (column F) with the new
#    "assertMatch(text,regex)" command inserted alphabetically.
#    It previously spanned F2:F44 (43 entries); it now spans F2:F45 (44).
# ---------------------------------------------------------------------------
$baseList = @(
    "appendText(var,appendWith)",
    "assertArrayContain(array,expected)",
    "assertArrayEqual(array1,array2,exactOrder)",
    "assertArrayNotContain(array,unexpected)",
    "assertContains(text,substring)",
    "assertCount(text,regex,expects)",
    "assertEmpty(text)",
    "assertEndsWith(text,suffix)",
    "assertEqual(expected,actual)",
    "assertMatch(text,regex)",
    "assertNotContain(text,substring)",
    "assertNotEmpty(text)",
    "assertNotEqual(expected,actual)",
    "assertStartsWith(text,prefix)",
    "assertTextOrder(var,descending)",
    "assertVarNotPresent(var)",
    "assertVarPresent(var)",
    "clear(vars)",
    "clearClipboard()",
    "copyFromClipboard(var)",
    "copyIntoClipboard(text)",
    "failImmediate(text)",
    "incrementChar(var,amount,config)",
    "macro(file,sheet,name)",
    "macroFlex(macro,input,output)",
    "outputToCloud(resource)",
    "prependText(var,prependWith)",
    "repeatUntil(steps,maxWaitMs)",
    "save(var,value)",
    "saveCount(text,regex,saveVar)",
    "saveMatches(text,regex,saveVar)",
    "saveReplace(text,regex,replace,saveVar)",
    "saveVariablesByPrefix(var,prefix)",
    "saveVariablesByRegex(var,regex)",
    "section(steps)",
    "split(text,delim,saveVar)",
    "startRecording()",
    "stopRecording()",
    "substringAfter(text,delim,saveVar)",
    "substringBefore(text,delim,saveVar)",
    "substringBetween(text,start,end,saveVar)",
    "verbose(text)",
    "waitFor(waitMs)",
    "waitForCondition(conditions,maxWaitMs)"
)
Set-ColumnValues $ws "F" 2 $baseList

# ---------------------------------------------------------------------------
# 4. Rewrite the "external" category list (column J) with the new
#    "openFile(filePath)" command inserted alphabetically.
#    It previously spanned J2:J6 (5 entries); it now spans J2:J7 (6).
# ---------------------------------------------------------------------------
$externalList = @(
    "openFile(filePath)",
    "runJUnit(className)",
    "runProgram(programPathAndParams)",
    "runProgramNoWait(programPathAndParams)",
    "tail(id,file)",
    "terminate(programName)"
)
Set-ColumnValues $ws "J" 2 $externalList

# ---------------------------------------------------------------------------
# 5. Update the named ranges that changed extents/columns.
# ---------------------------------------------------------------------------
$wb.Names.Item("base").RefersTo        = "='#system'!`$F`$2:`$F`$45"
$wb.Names.Item("external").RefersTo    = "='#system'!`$J`$2:`$J`$7"
$wb.Names.Item("target").RefersTo      = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo         = "='#system'!`$AA`$2:`$AA`$151"
$wb.Names.Item("webalert").RefersTo    = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo   = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo          = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo    = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo         = "='#system'!`$AF`$2:`$AF`$27"
